$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear rows 2 and 3 entirely (A2:AD3) so Excel omits them from sheetData on save
$ws.Range("A2:AD3").Clear()

# Update row 4 values in place (row number stays 4, values updated to match paper)
$ws.Range("A4").Value = 0.3299583781535884
$ws.Range("B4").Value = 0.3934589318997483
$ws.Range("C4").Value = 0.2830644034726374
$ws.Range("D4").Value = 0.3572424448617109
$ws.Range("E4").Value = 0.3017340007039827
$ws.Range("F4").Value = 7.009341074576589
$ws.Range("G4").Value = 9.144184282807482
$ws.Range("H4").Value = 5.437165512591095
$ws.Range("I4").Value = 7.920617889073113
$ws.Range("J4").Value = 6.060652977894222
$ws.Range("K4").Value = 2.312775817796077
$ws.Range("L4").Value = 3.160239842787972
$ws.Range("M4").Value = 1.702839254565729
$ws.Range("N4").Value = 2.665294820428289
$ws.Range("O4").Value = 1.944296884611536
$ws.Range("P4").Value = 65.12430999999999
$ws.Range("Q4").Value = 90.9821292620059
$ws.Range("R4").Value = 47.49237712278485
$ws.Range("S4").Value = 74.48617558543192
$ws.Range("T4").Value = 54.84334922681201
$ws.Range("U4").Value = 0.2013988016168057
$ws.Range("V4").Value = 0.3171759848229364
$ws.Range("W4").Value = 0.09930305172662894
$ws.Range("X4").Value = 0.2602628764248707
$ws.Range("Y4").Value = 0.1427292029185049
$ws.Range("Z4").Value = 0.8479314086336348
$ws.Range("AA4").Value = 0.9545603526709378
$ws.Range("AB4").Value = 0.6976235327873144
$ws.Range("AC4").Value = 0.9097342760093722
$ws.Range("AD4").Value = 0.7727318940014264
